$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 5.5
$ws.Range("N3").Value = 3.1
$ws.Range("O3").Value = 1.36
$ws.Range("R3").Value = 2.63
$ws.Range("S3").Value = 1.44
$ws.Range("U3").Value = 7
$ws.Range("V3").Value = 10
$ws.Range("AD3").Value = 9
$ws.Range("AE3").Value = 23
$ws.Range("AI3").Value = 67
$ws.Range("G5").Value = 2.47
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.65
$ws.Range("L5").Value = 1.27
$ws.Range("M5").Value = 3.1
$ws.Range("N5").Value = 1.8
$ws.Range("O5").Value = 1.8
$ws.Range("P5").Value = 1.39
$ws.Range("Q5").Value = 2.57
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 8.5
$ws.Range("W5").Value = 27
$ws.Range("X5").Value = 20
$ws.Range("Y5").Value = 28
$ws.Range("Z5").Value = 10.25
$ws.Range("AA5").Value = 6.3
$ws.Range("AB5").Value = 13
$ws.Range("AC5").Value = 55
$ws.Range("AD5").Value = 9.25
$ws.Range("AE5").Value = 14
$ws.Range("AF5").Value = 9.75
$ws.Range("AG5").Value = 30
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 28
$ws.Range("AJ5").Value = 400
$ws.Range("G6").Value = 1.65
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 4
$ws.Range("N6").Value = 1.29
$ws.Range("O6").Value = 3
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 2.8
$ws.Range("U6").Value = 13
$ws.Range("W6").Value = 17
$ws.Range("X6").Value = 11.5
$ws.Range("Y6").Value = 14.5
$ws.Range("Z6").Value = 26
$ws.Range("AA6").Value = 10.25
$ws.Range("AB6").Value = 11.25
$ws.Range("AC6").Value = 27
$ws.Range("AD6").Value = 26
$ws.Range("AE6").Value = 37
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 23
$ws.Range("G7").Value = 2.67
$ws.Range("I7").Value = 2.35
$ws.Range("T7").Value = 10
$ws.Range("U7").Value = 14.5
$ws.Range("V7").Value = 10
$ws.Range("W7").Value = 30
$ws.Range("X7").Value = 21
$ws.Range("Y7").Value = 27
$ws.Range("Z7").Value = 11.75
$ws.Range("AA7").Value = 6.7
$ws.Range("AD7").Value = 9
$ws.Range("AE7").Value = 12.5
$ws.Range("AF7").Value = 9.25
$ws.Range("AG7").Value = 24
$ws.Range("AH7").Value = 18
$ws.Range("AI7").Value = 26
$ws.Range("G14").Value = 3.1
$ws.Range("H14").Value = 3.35
$ws.Range("I14").Value = 2.1
$ws.Range("J14").Value = 1.05
$ws.Range("K14").Value = 7.6
$ws.Range("L14").Value = 1.27
$ws.Range("M14").Value = 3.4
$ws.Range("N14").Value = 1.8
$ws.Range("O14").Value = 1.9
$ws.Range("P14").Value = 1.39
$ws.Range("Q14").Value = 2.75
$ws.Range("T14").Value = 10.5
$ws.Range("U14").Value = 17.5
$ws.Range("V14").Value = 11
$ws.Range("Z14").Value = 7.6
$ws.Range("AA14").Value = 6.6
$ws.Range("AE14").Value = 10.75
$ws.Range("AG14").Value = 20
$ws.Range("AH14").Value = 16.5